# Edit: fix the method name shown in the "set/get" textbox on the
# "Object-Oriented Programming" slide (sldId 658) from "getGPA" to
# "getSalary". Only the text of the run containing "getGPA" is touched;
# every other run / formatting attribute on the shape is left intact.

$p = $ppt.ActivePresentation

$targetOld = "getGPA"
$targetNew = "getSalary"

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $shapes = $slide.Shapes

    for ($shi = 1; $shi -le $shapes.Count; $shi++) {
        $shape = $shapes.Item($shi)

        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $fullText = $tr.Text

            $idx = $fullText.IndexOf($targetOld)
            if ($idx -ge 0) {
                # Characters() is 1-based; replace only the matched run of
                # characters so surrounding runs/formatting are untouched.
                $chars = $tr.Characters($idx + 1, $targetOld.Length)
                $chars.Text = $targetNew
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not find text '$targetOld' in any slide to replace."
}
